$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 14 data - test case for APICORP (SuperNational Corporation)
$ws.Range("A14").Value = "XS2166383799"
$ws.Range("B14").Value = "APICORP (SuperNational Corporation)"
$ws.Range("D14").Value = "SA"
# Leading apostrophe forces a text/quote-prefixed cell (matches quotePrefix style in target)
$ws.Range("E14").Formula = "'Fixed Income, Government / Municipal"
$ws.Range("F14").Value = "Supernational in Saudi Arabia, use country of domicile as proxy i.e. SA (or Saudi Arabia)"

# Copy the highlighted "comment" style used elsewhere in column F onto the new note cell
$ws.Range("F9").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# Update the selected cell to match the saved workbook view
$ws.Range("B16").Select() | Out-Null
